# feat: add 2022-Q4 data
#
# 1. Insert a brand-new worksheet named "2022-Q4" right after "总计",
#    duplicating the layout/styling of the existing per-quarter sheet
#    (copy-then-overwrite keeps header/row styles identical), and fill
#    it with the two new fund rows.
# 2. Update the "总计" summary sheet: the 2022-Q4 totals become the new
#    row 2, and the previously-existing rows (2022-Q3 / 2022-Q1 / 2021-Q3)
#    each shift down by one row, with a brand-new row 5 created for the
#    2021-Q3 figures that used to live in row 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q4" worksheet right after "总计" by copying the
# sheet that currently sits in that slot ("2022-Q3") so the new sheet
# inherits identical column layout / header styling / page setup.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy([Type]::Missing, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Row 1 headers are already correct (copied). Overwrite the two data rows.
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "016174"
$q4Sheet.Range("C2").Value = "汇丰晋信策略优选混合A"
$q4Sheet.Range("D2").NumberFormat = "@"
$q4Sheet.Range("D2").Value = "1.78"
$q4Sheet.Range("E2").NumberFormat = "@"
$q4Sheet.Range("E2").Value = "74.92"
$q4Sheet.Range("F2").NumberFormat = "@"
$q4Sheet.Range("F2").Value = "2.19"
$q4Sheet.Range("G2").NumberFormat = "@"
$q4Sheet.Range("G2").Value = "0.0390"
$q4Sheet.Range("H2").Value = 10

# Row 3 is new - clone formatting from row 2, then set its values.
$q4Sheet.Range("A2:H2").Copy()
$q4Sheet.Range("A3:H3").PasteSpecial(-4122)  # xlPasteFormats

$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").NumberFormat = "@"
$q4Sheet.Range("B3").Value = "016175"
$q4Sheet.Range("C3").Value = "汇丰晋信策略优选混合C"
$q4Sheet.Range("D3").NumberFormat = "@"
$q4Sheet.Range("D3").Value = "0.40"
$q4Sheet.Range("E3").NumberFormat = "@"
$q4Sheet.Range("E3").Value = "74.92"
$q4Sheet.Range("F3").NumberFormat = "@"
$q4Sheet.Range("F3").Value = "2.19"
$q4Sheet.Range("G3").NumberFormat = "@"
$q4Sheet.Range("G3").Value = "0.0088"
$q4Sheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# Step 2: update "总计" - push existing rows down and add the 2022-Q4
# row at the top (row 2); a new row 5 is needed for 2021-Q3.
# ---------------------------------------------------------------------

# New row 5: clone style from row 4 first, then set 2021-Q3 figures
# (this used to be row 4's content).
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q3"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.13

# Row 4 becomes what used to be row 3 (2022-Q1).
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

# Row 3 becomes what used to be row 2 (2022-Q3).
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.17

# Row 2 becomes the brand-new 2022-Q4 totals.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.05
